$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.216.03'
$ws.Range('E2').Value = '  +1.71%  '
$ws.Range('D3').Value = '3.771.54'
$ws.Range('E3').Value = '  -0.52%  '
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').Value = "'624.10"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.88%  '
$ws.Range('D6').Value = "'163.90"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.56%  '
$ws.Range('D7').Value = '3.769.60'
$ws.Range('E7').Value = '  -0.51%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = "'0.520"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.63%  '
$ws.Range('D10').Value = "'0.160"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.18%  '
$ws.Range('D11').Value = "'0.451"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.32%  '
$ws.Range('D12').Value = "'6.63"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.72%  '
$ws.Range('D13').Value = "'0.0000247"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.19%  '
$ws.Range('D14').Value = "'35.37"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.66%  '
$ws.Range('D15').Value = '4.413.11'
$ws.Range('E15').Value = '  -0.33%  '
$ws.Range('D16').Value = '3.799.44'
$ws.Range('E16').Value = '  -1.00%  '
$ws.Range('D17').Value = '69.254.90'
$ws.Range('E17').Value = '  +1.77%  '
$ws.Range('D18').Value = "'17.78"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.78%  '
$ws.Range('D19').Value = "'7.08"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.36%  '
$ws.Range('E20').Value = '  -1.24%  '
$ws.Range('D21').Value = "'467.87"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.50%  '
$ws.Range('D22').Value = "'9.60"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.80%  '
$ws.Range('D23').Value = "'0.701"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').Value = "'0.0000149"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.18%  '
$ws.Range('D25').Value = "'83.01"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').Value = "'12.03"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.50%  '
$ws.Range('D27').Value = "'2.15"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.36%  '
$ws.Range('E28').Value = '  -0.07%  '
$ws.Range('D29').Value = "'9.97"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('D30').Value = '3.925.76'
$ws.Range('E30').Value = '  -0.42%  '
$ws.Range('D31').Value = "'2.66"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.15%  '
$ws.Range('D32').Value = "'2.23"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.09%  '
$ws.Range('D33').Value = "'7.27"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.75%  '
$ws.Range('D34').Value = "'28.85"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.29%  '
$ws.Range('D35').Value = "'0.998"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.22%  '
$ws.Range('B36').Value = 'RenzoRestakedETH'
$ws.Range('C36').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D36').Value = '3.726.91'
$ws.Range('E36').Value = '  -0.36%  '
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').Value = "'8.98"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.24%  '
$ws.Range('E38').Value = '  +2.23%  '
$ws.Range('D39').Value = "'0.152"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +9.17%  '
$ws.Range('D40').Value = "'3.33"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.80%  '
$ws.Range('D41').Value = "'5.79"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.43%  '
$ws.Range('E42').Value = '  -1.98%  '
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('D45').Value = "'0.300"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').Value = "'153.31"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.86%  '
$ws.Range('B47').Value = 'Arweave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D47').Value = "'43.09"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.19%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').Value = "'1.92"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.84%  '
$ws.Range('D49').Value = "'46.69"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.74%  '
$ws.Range('D50').Value = "'8.40"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.67%  '
$ws.Range('D51').Value = "'1.37"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.26%  '
